# Natmi following Dr Hou advice
#
# The sending/target cluster set grows from {FAPs, sCs} to {ECs, FAPs, sCs},
# so the 2x2 Slit3->Robo1 L-R pair table (rows 2-5) becomes a 3x3 table
# (rows 2-10). Ligand symbol is always "Slit3" and Receptor symbol is
# always "Robo1" for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ligand   = "Slit3"
$receptor = "Robo1"

# Sending-cluster block (columns E-J) - depends only on the sending cluster.
$sendingStats = @{
    "ECs"  = @(2, 0.6666666666666666, 2.096169666666666, 6.288508999999999, 0.01755614486036394, 0.01755614486036394)
    "FAPs" = @(3, 1, 108.7065656666667, 326.119697, 0.9104550287436967, 0.9104550287436968)
    "sCs"  = @(3, 1, 8.595326333333334, 25.785979, 0.07198882639593941, 0.07198882639593941)
}

# Target-cluster block (columns K-P) - depends only on the target cluster.
$targetStats = @{
    "ECs"  = @(2, 0.6666666666666666, 0.1325553333333333, 0.397666, 0.00533964316398423, 0.00533964316398423)
    "FAPs" = @(3, 1, 17.178266, 51.534798, 0.6919812904497951, 0.691981290449795)
    "sCs"  = @(3, 1, 7.513933666666667, 22.541801, 0.3026790663862208, 0.3026790663862208)
}

# Edge block (columns Q-T) - depends on the (sending cluster, target cluster) pair.
$edgeStats = @{
    "ECs|ECs"   = @(0.2778584688882222, 2.500726219994, 0.00009374354888955917, 0.00009374354888955917)
    "ECs|FAPs"  = @(36.00856011513133, 324.077041036182, 0.01214852377579818, 0.01214852377579817)
    "ECs|sCs"   = @(15.75047982941211, 141.754318464709, 0.005313877535676205, 0.005313877535676205)
    "FAPs|ECs"  = @(14.40963504746689, 129.686715427202, 0.004861504970346345, 0.004861504970346346)
    "FAPs|FAPs" = @(1867.390300968467, 16806.5127087162, 0.6300178456865685, 0.6300178456865685)
    "FAPs|sCs"  = @(816.8139235504774, 7351.325311954296, 0.2755756780867819, 0.275575678086782)
    "sCs|ECs"   = @(1.139356347223778, 10.254207125014, 0.0003843946447483253, 0.0003843946447483253)
    "sCs|FAPs"  = @(147.6528021108047, 1328.875218997242, 0.04981492098742842, 0.04981492098742842)
    "sCs|sCs"   = @(64.5847119120199, 581.262407208179, 0.02178951076376267, 0.02178951076376267)
}

$clusters = @("ECs", "FAPs", "sCs")

$row = 2
foreach ($sending in $clusters) {
    foreach ($target in $clusters) {
        $ws.Range("A$row").Value = $sending
        $ws.Range("B$row").Value = $ligand
        $ws.Range("C$row").Value = $receptor
        $ws.Range("D$row").Value = $target

        $s = $sendingStats[$sending]
        $ws.Range("E$row").Value = $s[0]
        $ws.Range("F$row").Value = $s[1]
        $ws.Range("G$row").Value = $s[2]
        $ws.Range("H$row").Value = $s[3]
        $ws.Range("I$row").Value = $s[4]
        $ws.Range("J$row").Value = $s[5]

        $t = $targetStats[$target]
        $ws.Range("K$row").Value = $t[0]
        $ws.Range("L$row").Value = $t[1]
        $ws.Range("M$row").Value = $t[2]
        $ws.Range("N$row").Value = $t[3]
        $ws.Range("O$row").Value = $t[4]
        $ws.Range("P$row").Value = $t[5]

        $e = $edgeStats["$sending|$target"]
        $ws.Range("Q$row").Value = $e[0]
        $ws.Range("R$row").Value = $e[1]
        $ws.Range("S$row").Value = $e[2]
        $ws.Range("T$row").Value = $e[3]

        $row = $row + 1
    }
}
